$wb = $excel.ActiveWorkbook

# Create the new sheet "2025-08-04" at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2025-08-04"

# Header row
$ws.Cells.Item(1, 1).Value = 'rank'
$ws.Cells.Item(1, 2).Value = 'title'
$ws.Cells.Item(1, 3).Value = 'author'
$ws.Cells.Item(1, 4).Value = 'latest_episode'

# Match header styling (bold, thin border, center/top align) used by the other daily sheets
$srcHeader = $wb.Worksheets.Item("2025-08-03").Range("A1:D1")
$srcHeader.Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Ranking data rows
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = '願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜'
$ws.Cells.Item(2, 3).Value = 'ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)'
$ws.Cells.Item(2, 4).Value = '第5話-1：魔法のお勉強'
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = '姫様“拷問”の時間です'
$ws.Cells.Item(3, 3).Value = '原作:春原ロビンソン　漫画:ひらけい'
$ws.Cells.Item(3, 4).Value = '拷問146'
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = '遊遊じてき。'
$ws.Cells.Item(4, 3).Value = 'カンケー'
$ws.Cells.Item(4, 4).Value = '第2話'
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = '我輩は猫魔導師である ～キジトラ・ルークの快適チート猫生活～'
$ws.Cells.Item(5, 3).Value = '原作：猫神信仰研究会 漫画：三國大和 キャラクター原案：ハム'
$ws.Cells.Item(5, 4).Value = '第23話'
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'いとこのこ'
$ws.Cells.Item(6, 3).Value = 'いぬちく(著者)'
$ws.Cells.Item(6, 4).Value = '第36話'
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = '窓際編集とバカにされた俺が、双子ＪＫと同居することになった'
$ws.Cells.Item(7, 3).Value = 'うさおとめ(著者) 茨木野(原作) トモゼロ(キャラクター原案)'
$ws.Cells.Item(7, 4).Value = '第5話①'
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = '不純な彼女達は懺悔しない'
$ws.Cells.Item(8, 3).Value = 'ポロロッカ(著者)'
$ws.Cells.Item(8, 4).Value = '第29話'
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = '悪役貴族として必要なそれ'
$ws.Cells.Item(9, 3).Value = 'まさこりん(原作) 夏野うみ(作画) 村カルキ(キャラクターデザイン)'
$ws.Cells.Item(9, 4).Value = '第18話①'
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = '転生したらスライムだった件　クレイマンREVENGE'
$ws.Cells.Item(10, 3).Value = '原作：伏瀬 漫画：カジカ航 キャラクター原案：みっつばー'
$ws.Cells.Item(10, 4).Value = '第28話　魔人形'
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'おんなのこのけんをてにいれた'
$ws.Cells.Item(11, 3).Value = '福岡太朗(著者)'
$ws.Cells.Item(11, 4).Value = '16本目'
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'え、社内システム全てワンオペしている私を解雇ですか？'
$ws.Cells.Item(12, 3).Value = '漫画：伊於 原作：下城米雪 キャラクター原案：icchi'
$ws.Cells.Item(12, 4).Value = '4巻発売告知漫画'
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版'
$ws.Cells.Item(13, 3).Value = '漫画/すたひろ 原作/Y.A'
$ws.Cells.Item(13, 4).Value = 'chapter66【35話①】'
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'まったく最近の探偵ときたら'
$ws.Cells.Item(14, 3).Value = '五十嵐正邦(著者)'
$ws.Cells.Item(14, 4).Value = '第114話'
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'ギャルとダンジョンと周回遅れの探索英雄譚'
$ws.Cells.Item(15, 3).Value = '漫画家： 水田ケンジ 原作：榊一郎 キャラクター原案：黒獅子'
$ws.Cells.Item(15, 4).Value = '第1話'
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'リビルドワールド'
$ws.Cells.Item(16, 3).Value = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$ws.Cells.Item(16, 4).Value = '第71話④'
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = '地味子な三葉さんが僕を誘惑する'
$ws.Cells.Item(17, 3).Value = 'はぶらえる(著者)'
$ws.Cells.Item(17, 4).Value = '第10話後半'
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = '水魔法なんて使えないと追放されたけど、水が万能だと気がつき水の賢者と呼ばれるまでに成長しました'
$ws.Cells.Item(18, 3).Value = '原作：空地大乃 漫画：結ゆい キャラクター原案：神吉李花'
$ws.Cells.Item(18, 4).Value = '第8話'
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = '異世界迷宮のオーパーツ'
$ws.Cells.Item(19, 3).Value = '三狛ハル(著者)'
$ws.Cells.Item(19, 4).Value = '第3話-①：おそらく高貴な布'
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = '女友達は頼めば意外とヤらせてくれる'
$ws.Cells.Item(20, 3).Value = 'ろくろ(漫画) 鏡遊(原作)'
$ws.Cells.Item(20, 4).Value = '第23話②'
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = '魔都精兵のスレイブ'
$ws.Cells.Item(21, 3).Value = '原作:タカヒロ　漫画:竹村洋平'
$ws.Cells.Item(21, 4).Value = '第158話　黄泉津大山'
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = 'Ｓ級ギルドを追放されたけど、実は俺だけドラゴンの言葉がわかるので、気付いたときには竜騎士の頂点を極めてました。'
$ws.Cells.Item(22, 3).Value = 'ひそな(漫画) 三木なずな(原作) 白狼(キャラクター原案)'
$ws.Cells.Item(22, 4).Value = '第37話-2'
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = '半人前の恋人'
$ws.Cells.Item(23, 3).Value = '川田大智'
$ws.Cells.Item(23, 4).Value = '第49話'
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = '宮廷鍛冶師の幸せな日常 ～ブラックな職場を追放されたが、隣国で公爵令嬢に溺愛されながらホワイトな生活送ります～'
$ws.Cells.Item(24, 3).Value = '上林眞(著者) 木嶋隆太(原作) a20(キャラクター原案)'
$ws.Cells.Item(24, 4).Value = '第33話-②'
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = '魔のものたちは企てる'
$ws.Cells.Item(25, 3).Value = '加藤拓弐(原作) ガしガし(作画)'
$ws.Cells.Item(25, 4).Value = '第28話'
$ws.Cells.Item(26, 1).Value = 25
$ws.Cells.Item(26, 2).Value = '無敵商人の異世界成り上がり物語 ～現代の製品を自在に取り寄せるスキルがあるので異世界では楽勝です～'
$ws.Cells.Item(26, 3).Value = '隆原ヒロタ(漫画) 青山有(原作) ぷきゅのすけ(キャラクターデザイン)'
$ws.Cells.Item(26, 4).Value = '第35話②'
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = '愚かな天使は悪魔と踊る'
$ws.Cells.Item(27, 3).Value = 'アズマサワヨシ(著者)'
$ws.Cells.Item(27, 4).Value = '第100話④'
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 'リアリスト魔王による聖域なき異世界改革'
$ws.Cells.Item(28, 3).Value = '鈴木マナツ(漫画) 羽田遼亮(原作) ゆーげん(キャラクターデザイン) ひたきゆう(キャラクターデザイン)'
$ws.Cells.Item(28, 4).Value = '第67幕④'
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = '世界の終わりの世界録(アンコール)'
$ws.Cells.Item(29, 3).Value = '雨水龍(著者) 細音啓(原作) ふゆの春秋(キャラクター原案)'
$ws.Cells.Item(29, 4).Value = '第96話①'
$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = '十歳の最強魔導師'
$ws.Cells.Item(30, 3).Value = '猫月 天乃聖樹'
$ws.Cells.Item(30, 4).Value = '第1話'
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = '育成上手な冒険者、幼女を拾い、セカンドライフを育児に捧げる'
$ws.Cells.Item(31, 3).Value = '原作／リズ 漫画／森見明日'
$ws.Cells.Item(31, 4).Value = '第15話'
$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = 'アラサーがVTuberになった話。'
$ws.Cells.Item(32, 3).Value = '犬威赤彦(漫画) とくめい(原作) カラスBTK(キャラクター原案)'
$ws.Cells.Item(32, 4).Value = '第25話'
$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = '俺堕ちスレイブヒーローコレクション'
$ws.Cells.Item(33, 3).Value = 'ゆっ栗栖(著者)'
$ws.Cells.Item(33, 4).Value = '第11話後半'
$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = '追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。'
$ws.Cells.Item(34, 3).Value = '六志麻あさ 業務用餅 kisui'
$ws.Cells.Item(34, 4).Value = '第６９話'
$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = '貴方は猫（わたし）の下僕です ～ねことげぼくのヒミツなカンケイ～'
$ws.Cells.Item(35, 3).Value = '大田優一(著者)'
$ws.Cells.Item(35, 4).Value = '第14話後半'
$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = '転生してあらゆるモノに好かれながら異世界で好きな事をして生きて行く'
$ws.Cells.Item(36, 3).Value = '都尾琉(漫画) 御峰。(原作)'
$ws.Cells.Item(36, 4).Value = '第27話①'
$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = '二番目な僕と一番の彼女'
$ws.Cells.Item(37, 3).Value = 'ぬずタニ(漫画) 和尚(原作) ミュシャ(キャラクター原案)'
$ws.Cells.Item(37, 4).Value = '第1話'
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = '経験値貯蓄でのんびり傷心旅行 ～勇者と恋人に追放された戦士の無自覚ざまぁ～'
$ws.Cells.Item(38, 3).Value = '奏ヨシキ(著者) 徳川レモン(原作) riritto(キャラクターデザイン)'
$ws.Cells.Item(38, 4).Value = '第37話-2'
$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = '小林さんちのメイドラゴン'
$ws.Cells.Item(39, 3).Value = 'クール教信者'
$ws.Cells.Item(39, 4).Value = '第147話'
$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = 'チュートリアルが始まる前に ボスキャラ達を破滅させない為に俺ができる幾つかの事'
$ws.Cells.Item(40, 3).Value = '横山コウヂ(漫画) 髙橋炬燵(原作) カカオ・ランタン(キャラクターデザイン)'
$ws.Cells.Item(40, 4).Value = '第14話①'
$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 'ハーレムより平穏を！異世界で静かにニート姫させてくれ'
$ws.Cells.Item(41, 3).Value = 'さかたはるき(原作) かわやばぐ(作画)'
$ws.Cells.Item(41, 4).Value = '第14話前半'
$ws.Cells.Item(42, 1).Value = 41
$ws.Cells.Item(42, 2).Value = '理想のヒモ生活'
$ws.Cells.Item(42, 3).Value = '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)'
$ws.Cells.Item(42, 4).Value = '第86話　その2'
$ws.Cells.Item(43, 1).Value = 42
$ws.Cells.Item(43, 2).Value = 'くらいあの子としたいこと'
$ws.Cells.Item(43, 3).Value = '碇マナツ(著者)'
$ws.Cells.Item(43, 4).Value = '第81話'
$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = '底辺ハンターが【リターン】スキルで現代最強 ～前世の知識と死に戻りを駆使して、人類最速レベルアップ～'
$ws.Cells.Item(44, 3).Value = '原作：萩鵜アキ 漫画：仲間友 キャラクター原案：gunkan'
$ws.Cells.Item(44, 4).Value = '第19話'
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$ws.Cells.Item(45, 3).Value = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$ws.Cells.Item(45, 4).Value = '第81話その1'
$ws.Cells.Item(46, 1).Value = 45
$ws.Cells.Item(46, 2).Value = 'きみの願いが叶うまで'
$ws.Cells.Item(46, 3).Value = '浅月のりと(著者)'
$ws.Cells.Item(46, 4).Value = '第4話-1'
$ws.Cells.Item(47, 1).Value = 46
$ws.Cells.Item(47, 2).Value = '花子さん、食べないで'
$ws.Cells.Item(47, 3).Value = '茸谷きの子(著者)'
$ws.Cells.Item(47, 4).Value = '第9話前半'
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = '最強勇者パーティーは愛が知りたい'
$ws.Cells.Item(48, 3).Value = '山田肌襦袢'
$ws.Cells.Item(48, 4).Value = '第28話「最後はこぶしがあればいい」'
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = '豚のレバーは加熱しろ'
$ws.Cells.Item(49, 3).Value = 'みなみ(漫画) 逆井卓馬(原作) 遠坂あさぎ(キャラクターデザイン)'
$ws.Cells.Item(49, 4).Value = '第42話②'
$ws.Cells.Item(50, 1).Value = 49
$ws.Cells.Item(50, 2).Value = 'ホームセンターごと呼び出された私の大迷宮リノベーション！'
$ws.Cells.Item(50, 3).Value = 'ばたっち(漫画) 星崎崑(原作) 志田(キャラクター原案)'
$ws.Cells.Item(50, 4).Value = '番外編①'
$ws.Cells.Item(51, 1).Value = 50
$ws.Cells.Item(51, 2).Value = '角と板と魔法記師 Tab.4'
$ws.Cells.Item(51, 3).Value = 'とりから'
$ws.Cells.Item(51, 4).Value = '第25話の6'

# Restore the originally active sheet/selection so the workbook-level view state is unaffected
$wb.Worksheets.Item("Sheet1").Activate()
